$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of row 2 (A2:B2) without shifting row 3 up
$ws.Range("A2:B2").ClearContents()

# Update cell B3 text value from "1" to "1.0" (kept as text, not a number)
$ws.Range("B3").Value = "'1.0"
$ws.Range("B3").ClearFormats()
